$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 "Marking": Right column 5 -> 4, Wrong column -1 -> -2
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 "Total": Right column 95 -> 76, Wrong column -2 -> -4
$ws.Range("B12").Value = 76
$ws.Range("C12").Value = -4

# Update the "Max" column total text to reflect corrected totals
$ws.Range("E12").Value = "72 / 112"
